# Fixar (re-sincronizar) a ordem das UFs nas linhas 21-27 (cauda do ranking)
# das abas tot-arrecad, avg-arrecad, max-arrecad e tx-sucesso, para o
# mesmo conjunto de dados de autoria/mencoes do ano de analise.

$wb = $excel.ActiveWorkbook

# tot-arrecad (sheet2): linhas 21-27
$ws = $wb.Worksheets.Item("tot-arrecad")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "TO"
$ws.Range("A23").Value = "AP"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# avg-arrecad (sheet3): linhas 21-27
$ws = $wb.Worksheets.Item("avg-arrecad")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "TO"
$ws.Range("A23").Value = "AP"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"
$ws.Range("A26").Value = "RN"
$ws.Range("A27").Value = "PI"

# max-arrecad (sheet4): linhas 21-25 (26-27 ja corretas)
$ws = $wb.Worksheets.Item("max-arrecad")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "TO"
$ws.Range("A23").Value = "AP"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"

# tx-sucesso (sheet5): linhas 21-27
$ws = $wb.Worksheets.Item("tx-sucesso")
$ws.Range("A21").Value = "RO"
$ws.Range("A22").Value = "TO"
$ws.Range("A23").Value = "AP"
$ws.Range("A24").Value = "AL"
$ws.Range("A25").Value = "MT"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"
